# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.967.17"
$ws.Range("E2").Value = "  -3.98%  "

$ws.Range("D3").Value = "'3.295.65"
$ws.Range("E3").Value = "  -5.78%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'541.59"
$ws.Range("E5").Value = "  -2.34%  "

$ws.Range("D6").Value = "'170.65"
$ws.Range("E6").Value = "  -4.64%  "

$ws.Range("D7").Value = "'0.611"
$ws.Range("E7").Value = "  -4.06%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").Value = "'3.281.40"
$ws.Range("E9").Value = "  -6.00%  "

$ws.Range("D10").Value = "'0.608"
$ws.Range("E10").Value = "  -3.97%  "

$ws.Range("D11").Value = "'0.151"
$ws.Range("E11").Value = "  -0.79%  "

$ws.Range("D12").Value = "'52.57"
$ws.Range("E12").Value = "  -2.28%  "

$ws.Range("D13").Value = "'0.0000263"
$ws.Range("E13").Value = "  -2.78%  "

$ws.Range("D14").Value = "'8.80"
$ws.Range("E14").Value = "  -4.79%  "

$ws.Range("D15").Value = "'3.796.73"
$ws.Range("E15").Value = "  -6.31%  "

$ws.Range("D16").Value = "'17.89"
$ws.Range("E16").Value = "  -3.07%  "

$ws.Range("E17").Value = "  -4.13%  "

$ws.Range("D18").Value = "'3.274.08"
$ws.Range("E18").Value = "  -6.19%  "

$ws.Range("D19").Value = "'11.57"
$ws.Range("E19").Value = "  -4.29%  "

$ws.Range("D20").Value = "'62.839.09"
$ws.Range("E20").Value = "  -4.20%  "

$ws.Range("D21").Value = "'0.965"
$ws.Range("E21").Value = "  -2.79%  "

$ws.Range("D22").Value = "'411.28"
$ws.Range("E22").Value = "  -1.69%  "

$ws.Range("D23").Value = "'4.38"
$ws.Range("E23").Value = "  +6.23%  "

$ws.Range("D24").Value = "'4.00"
$ws.Range("E24").Value = "  -1.18%  "

$ws.Range("D25").Value = "'13.48"
$ws.Range("E25").Value = "  +4.79%  "

$ws.Range("D26").Value = "'82.40"
$ws.Range("E26").Value = "  -4.75%  "

$ws.Range("D27").Value = "'10.49"
$ws.Range("E27").Value = "  -2.84%  "

$ws.Range("D28").Value = "'2.70"
$ws.Range("E28").Value = "  -5.47%  "

$ws.Range("D29").Value = "'8.54"
$ws.Range("E29").Value = "  -5.91%  "

$ws.Range("D30").Value = "'28.87"
$ws.Range("E30").Value = "  -4.76%  "

$ws.Range("D31").Value = "'6.30"
$ws.Range("E31").Value = "  -3.41%  "

$ws.Range("D32").Value = "'11.24"
$ws.Range("E32").Value = "  -4.18%  "

$ws.Range("D33").Value = "'572.93"
$ws.Range("E33").Value = "  -5.84%  "

$ws.Range("E34").Value = "  -4.00%  "

$ws.Range("D35").Value = "'57.49"
$ws.Range("E35").Value = "  -3.46%  "

$ws.Range("E36").Value = "  +0.17%  "

$ws.Range("D37").Value = "'0.145"
$ws.Range("E37").Value = "  -1.23%  "

$ws.Range("D38").Value = "'34.86"
$ws.Range("E38").Value = "  -6.93%  "

$ws.Range("D39").Value = "'3.36"
$ws.Range("E39").Value = "  +3.33%  "

$ws.Range("D40").Value = "'0.0₃0731"
$ws.Range("E40").Value = "  -7.50%  "

$ws.Range("D41").Value = "'0.363"
$ws.Range("E41").Value = "  -4.81%  "

$ws.Range("D42").Value = "'3.100.83"
$ws.Range("E42").Value = "  -8.60%  "

$ws.Range("E43").Value = "  -0.24%  "

$ws.Range("B44").Value = "ThetaToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D44").Value = "'2.74"
$ws.Range("E44").Value = "  -4.40%  "

$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "'3.22"
$ws.Range("E45").Value = "  -0.22%  "

$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "'2.40"
$ws.Range("E46").Value = "  -5.42%  "

$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "'0.0397"
$ws.Range("E47").Value = "  -4.25%  "

$ws.Range("E48").Value = "  -4.19%  "

$ws.Range("E49").Value = "  -3.60%  "

$ws.Range("D50").Value = "'132.02"
$ws.Range("E50").Value = "  -4.31%  "

$ws.Range("D51").Value = "'7.94"
$ws.Range("E51").Value = "  -6.46%  "
